$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column: DA corresponds to 2026-01-14 (serial 46036), the day after
# the existing last column CZ (2026-01-13, serial 46035).

# Header cell (row 1): date value, same formatting as CZ1 (vertical-centered
# date format). Value is set before formatting so dependent formulas
# (COUNTA/COUNTIF over the full row range) see the write and recalc.
$h = $ws.Range("DA1")
$h.Value = 46036
$h.NumberFormat = "mm-dd-yy"
$h.VerticalAlignment = -4108

# Attendance rows: default is "P" (Présent), matching each player's usual
# style (horizontally centered). A few players have a different status for
# this date.
$rowsB = @(13, 26)   # "B"
$rowsR = @(25)       # "R"
$rowsP = @(2,3,4,5,6,7,8,9,10,11,14,15,16,17,18,19,20,22,24,27,28,29,30)   # "P"

foreach ($r in $rowsP) {
    $c = $ws.Cells.Item($r, 105)
    $c.Value = "P"
    $c.HorizontalAlignment = -4108
}

foreach ($r in $rowsB) {
    $c = $ws.Cells.Item($r, 105)
    $c.Value = "B"
    $c.HorizontalAlignment = -4108
}

foreach ($r in $rowsR) {
    $c = $ws.Cells.Item($r, 105)
    $c.Value = "R"
    $c.HorizontalAlignment = -4108
}

# Match the author's final cursor position (row 22 of the new column).
$ws.Range("DD22").Select() | Out-Null
